$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '65.078.53'
$ws.Range('E2').Value = '  +0.53%  '
$ws.Range('D3').Value = '3.528.14'
$ws.Range('E3').Value = '  -0.43%  '
$ws.Range('E4').Value = '  +0.07%  '
Set-TextValue 'D5' '593.35'
$ws.Range('E5').Value = '  -0.87%  '
Set-TextValue 'D6' '134.13'
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('D7').Value = '3.527.24'
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('E8').Value = '  +0.10%  '
Set-TextValue 'D9' '0.492'
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('E10').Value = '  +2.09%  '
$ws.Range('E11').Value = '  +3.57%  '
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').Value = '4.127.36'
$ws.Range('E13').Value = '  -0.25%  '
Set-TextValue 'D14' '27.74'
$ws.Range('E14').Value = '  +2.96%  '
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('E16').Value = '  +0.61%  '
$ws.Range('D17').Value = '3.522.79'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').Value = '65.059.87'
$ws.Range('E18').Value = '  +0.66%  '
Set-TextValue 'D19' '10.18'
$ws.Range('E19').Value = '  +1.54%  '
Set-TextValue 'D20' '14.46'
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('E21').Value = '  -1.84%  '
Set-TextValue 'D22' '392.45'
$ws.Range('E22').Value = '  +1.65%  '
Set-TextValue 'D23' '0.582'
$ws.Range('E23').Value = '  +1.16%  '
Set-TextValue 'D24' '75.02'
$ws.Range('E24').Value = '  +1.18%  '
$ws.Range('D25').Value = '3.672.85'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  -2.47%  '
Set-TextValue 'D28' '7.76'
$ws.Range('E28').Value = '  +2.03%  '
$ws.Range('E29').Value = '  +11.41%  '
$ws.Range('E30').Value = '  -0.21%  '
$ws.Range('E31').Value = '  -0.30%  '
Set-TextValue 'D32' '8.34'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').Value = '3.538.67'
Set-TextValue 'D34' '24.09'
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +0.45%  '
$ws.Range('E37').Value = '  +6.57%  '
$ws.Range('E38').Value = '  +2.80%  '
Set-TextValue 'D40' '168.38'
$ws.Range('E40').Value = '  -0.77%  '
Set-TextValue 'D41' '0.0812'
$ws.Range('E41').Value = '  +1.21%  '
Set-TextValue 'D42' '0.823'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('E43').Value = '  +6.12%  '
Set-TextValue 'D44' '25.87'
$ws.Range('E44').Value = '  -3.35%  '
Set-TextValue 'D45' '42.99'
$ws.Range('E45').Value = '  +1.02%  '
Set-TextValue 'D46' '0.999'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('E48').Value = '  +1.84%  '
Set-TextValue 'D49' '6.91'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').Value = '2.425.85'
$ws.Range('E50').Value = '  -0.88%  '
$ws.Range('E51').Value = '  +6.60%  '
